$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 35021
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 35021
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 35021
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -35957
$ws.Range("H23").Value = 35021
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 35021
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 35021
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -35489
$ws.Range("H28").Value = 1966
$ws.Range("I28").Value = 1450
$ws.Range("J28").Value = 2998
$ws.Range("K28").Value = 1450
$ws.Range("L28").Value = 2998
$ws.Range("M28").Value = -965
$ws.Range("N28").Value = -3968
$ws.Range("H32").Value = 11499.75
$ws.Range("J32").Value = 12999.667
$ws.Range("L32").Value = 12999.667
$ws.Range("N32").Value = -13651.667
$ws.Range("H70").Value = 1016.8333
$ws.Range("I70").Value = 1025.25
$ws.Range("K70").Value = 3075.75
$ws.Range("M70").Value = -2805.75
$ws.Range("H73").Value = 1016.8333
$ws.Range("I73").Value = 1025.25
$ws.Range("K73").Value = 3075.75
$ws.Range("M73").Value = -2139.75
$ws.Range("H86").Value = 4985
$ws.Range("J86").Value = 4984
$ws.Range("L86").Value = 4984
$ws.Range("N86").Value = -7230
$ws.Range("H89").Value = 4985
$ws.Range("J89").Value = 4984
$ws.Range("L89").Value = 24920
$ws.Range("N89").Value = -36152
$ws.Range("H107").Value = 617.55554
$ws.Range("I107").Value = 676.6667
$ws.Range("K107").Value = 676.6667
$ws.Range("M107").Value = 1243.3333
$ws.Range("H141").Value = 7428.143
$ws.Range("I141").Value = 7999.5
$ws.Range("J141").Value = 7199.6
$ws.Range("K141").Value = 23998.5
$ws.Range("L141").Value = 21598.8
$ws.Range("M141").Value = -18818.5
$ws.Range("N141").Value = -31958.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8334.446
$ws.Range("I32").Value = 8080.028
$ws.Range("K32").Value = 8080.028
$ws.Range("M32").Value = -7793.028
$ws.Range("H132").Value = 2902.4707
$ws.Range("I132").Value = 2334.9
$ws.Range("K132").Value = 7004.700000000001
$ws.Range("M132").Value = -4474.700000000001
$ws.Range("H138").Value = 99890.664
$ws.Range("J138").Value = 99890.664
$ws.Range("L138").Value = 99890.664
$ws.Range("N138").Value = -110170.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 5380.6
$ws.Range("I11").Value = 252
$ws.Range("K11").Value = 252
$ws.Range("M11").Value = -112
$ws.Range("H20").Value = 9799.333000000001
$ws.Range("I20").Value = 9799.333000000001
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 9799.333000000001
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -9552.333000000001
$ws.Range("N20").ClearContents()
$ws.Range("H86").Value = 4409.875
$ws.Range("I86").Value = 4018.4285
$ws.Range("J86").Value = 7150
$ws.Range("K86").Value = 4018.4285
$ws.Range("L86").Value = 7150
$ws.Range("M86").Value = -2895.4285
$ws.Range("N86").Value = -9396
$ws.Range("H89").Value = 4409.875
$ws.Range("I89").Value = 4018.4285
$ws.Range("J89").Value = 7150
$ws.Range("K89").Value = 20092.1425
$ws.Range("L89").Value = 35750
$ws.Range("M89").Value = -14476.1425
$ws.Range("N89").Value = -46982
$ws.Range("H94").Value = 954
$ws.Range("I94").Value = 614.1429000000001
$ws.Range("J94").Value = 3333
$ws.Range("K94").Value = 614.1429000000001
$ws.Range("L94").Value = 3333
$ws.Range("M94").Value = -163.1429000000001
$ws.Range("N94").Value = -4235
$ws.Range("H134").Value = 2378.9524
$ws.Range("I134").Value = 2331.111
$ws.Range("K134").Value = 6993.333
$ws.Range("M134").Value = -4458.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H22").Value = 732.8570999999999
$ws.Range("I22").Value = 688.3333
$ws.Range("K22").Value = 688.3333
$ws.Range("M22").Value = -338.3333
$ws.Range("H25").Value = 20013
$ws.Range("J25").Value = 20013
$ws.Range("L25").Value = 20013
$ws.Range("N25").Value = -20361
$ws.Range("H41").Value = 39643.332
$ws.Range("J41").Value = 39643.332
$ws.Range("L41").Value = 39643.332
$ws.Range("N41").Value = -40499.332
$ws.Range("H58").Value = 2912.4443
$ws.Range("I58").Value = 2883.4285
$ws.Range("K58").Value = 2883.4285
$ws.Range("M58").Value = -2680.4285
$ws.Range("H132").Value = 4513.2
$ws.Range("I132").Value = 4391.75
$ws.Range("K132").Value = 13175.25
$ws.Range("M132").Value = -10645.25
$ws.Range("H134").Value = 2968
$ws.Range("I134").Value = 3161.6
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 9484.799999999999
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -6949.799999999999
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 2912.4443
$ws.Range("I136").Value = 2883.4285
$ws.Range("K136").Value = 8650.2855
$ws.Range("M136").Value = -6100.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 89.8
$ws.Range("J2").Value = 99.75
$ws.Range("L2").Value = 598.5
$ws.Range("N2").Value = -824.5
$ws.Range("H55").Value = 1400
$ws.Range("I55").Value = 1400
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 4200
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -4023
$ws.Range("N55").ClearContents()
$ws.Range("H139").Value = 1046.6666
$ws.Range("I139").Value = 1046.6666
$ws.Range("K139").Value = 3139.9998
$ws.Range("M139").Value = 2000.0002
$ws.Range("H140").Value = 1314.9286
$ws.Range("I140").Value = 1314.9286
$ws.Range("K140").Value = 3944.7858
$ws.Range("M140").Value = 1235.2142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 2998.6667
$ws.Range("I102").Value = 2998.5
$ws.Range("J102").Value = 2999
$ws.Range("K102").Value = 2998.5
$ws.Range("L102").Value = 2999
$ws.Range("M102").Value = -1376.5
$ws.Range("N102").Value = -6243
$ws.Range("H113").Value = 2723.5
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 3497
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 3497
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -7837
$ws.Range("H126").Value = 2899
$ws.Range("J126").Value = 2899
$ws.Range("L126").Value = 8697
$ws.Range("N126").Value = -13637
$ws.Range("H132").Value = 3869
$ws.Range("I132").Value = 3680.8
$ws.Range("J132").Value = 4406.7144
$ws.Range("K132").Value = 11042.4
$ws.Range("L132").Value = 13220.1432
$ws.Range("M132").Value = -8512.400000000001
$ws.Range("N132").Value = -18280.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5763.2
$ws.Range("I22").Value = 4866.6665
$ws.Range("J22").Value = 6147.4287
$ws.Range("K22").Value = 4866.6665
$ws.Range("L22").Value = 6147.4287
$ws.Range("M22").Value = -4571.6665
$ws.Range("N22").Value = -6737.4287
$ws.Range("H27").Value = 5763.2
$ws.Range("I27").Value = 4866.6665
$ws.Range("J27").Value = 6147.4287
$ws.Range("K27").Value = 4866.6665
$ws.Range("L27").Value = 6147.4287
$ws.Range("M27").Value = -4759.6665
$ws.Range("N27").Value = -6361.4287
$ws.Range("H136").Value = 3839.9644
$ws.Range("I136").Value = 4066.1365
$ws.Range("J136").Value = 3010.6667
$ws.Range("K136").Value = 12198.4095
$ws.Range("L136").Value = 9032.000100000001
$ws.Range("M136").Value = -9648.4095
$ws.Range("N136").Value = -14132.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 67599.8
$ws.Range("I46").Value = 70000
$ws.Range("J46").Value = 66999.75
$ws.Range("K46").Value = 70000
$ws.Range("L46").Value = 66999.75
$ws.Range("M46").Value = -69769
$ws.Range("N46").Value = -67461.75
$ws.Range("H62").Value = 167262
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 167262
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 50734.145
$ws.Range("I70").Value = 47989
$ws.Range("J70").Value = 51191.668
$ws.Range("K70").Value = 47989
$ws.Range("L70").Value = 51191.668
$ws.Range("M70").Value = -47674
$ws.Range("N70").Value = -51821.668
$ws.Range("H73").Value = 50734.145
$ws.Range("I73").Value = 47989
$ws.Range("J73").Value = 51191.668
$ws.Range("K73").Value = 47989
$ws.Range("L73").Value = 51191.668
$ws.Range("M73").Value = -46897
$ws.Range("N73").Value = -53375.668
$ws.Range("H81").Value = 5387.9287
$ws.Range("I81").Value = 3119.75
$ws.Range("K81").Value = 6239.5
$ws.Range("M81").Value = -5178.5
$ws.Range("H84").Value = 5387.9287
$ws.Range("I84").Value = 3119.75
$ws.Range("K84").Value = 31197.5
$ws.Range("M84").Value = -25893.5
$ws.Range("H122").Value = 4740.6665
$ws.Range("I122").Value = 4611
$ws.Range("K122").Value = 13833
$ws.Range("M122").Value = -11383
$ws.Range("H126").Value = 2090.9333
$ws.Range("I126").Value = 2133.1428
$ws.Range("K126").Value = 6399.428400000001
$ws.Range("M126").Value = -3929.428400000001
$ws.Range("H134").Value = 67599.8
$ws.Range("I134").Value = 70000
$ws.Range("J134").Value = 66999.75
$ws.Range("K134").Value = 210000
$ws.Range("L134").Value = 200999.25
$ws.Range("M134").Value = -207465
$ws.Range("N134").Value = -206069.25
